$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ((Intercept))
$ws.Range("B2").Value = 1.518349
$ws.Range("D2").Value = 4.327757
$ws.Range("E2").Value = 0.040618

# Row 3 (household_group_collapsed)
$ws.Range("B3").Value = 15.268191
$ws.Range("D3").Value = 21.759488
$ws.Range("E3").Value = 0

# Row 4 (Residuals)
$ws.Range("B4").Value = 28.768867
$ws.Range("C4").Value = 82

# Row 5 (SM-Control)
$ws.Range("G5").Value = 0.08558
$ws.Range("H5").Value = -0.521763
$ws.Range("I5").Value = 0.692923
$ws.Range("J5").Value = 0.939591

# Row 6 (SM + Traps-Control)
$ws.Range("G6").Value = 0.945603
$ws.Range("H6").Value = 0.377631
$ws.Range("I6").Value = 1.513574
$ws.Range("J6").Value = 0.000441

# Row 7 (SM + Traps-SM)
$ws.Range("G7").Value = 0.860022
$ws.Range("H7").Value = 0.513164
$ws.Range("I7").Value = 1.206881
$ws.Range("J7").Value = 0
